# Apply 2023-10-25 violent crime data update
# Updates the 2023 (column J) running totals across the Citywide Totals,
# By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6292
$ws.Range("J3").Value = 6688
$ws.Range("J4").Value = 1449
$ws.Range("J5").Value = 514
$ws.Range("J6").Value = 8791
$ws.Range("J7").Value = 23734

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 188
$ws.Range("J4").Value = 106
$ws.Range("J8").Value = 1500
$ws.Range("J9").Value = 122
$ws.Range("J10").Value = 176
$ws.Range("J11").Value = 396
$ws.Range("J15").Value = 268
$ws.Range("J16").Value = 96
$ws.Range("J19").Value = 701
$ws.Range("J20").Value = 492
$ws.Range("J23").Value = 222
$ws.Range("J27").Value = 146
$ws.Range("J29").Value = 1302
$ws.Range("J33").Value = 1074
$ws.Range("J36").Value = 320
$ws.Range("J37").Value = 725
$ws.Range("J40").Value = 52
$ws.Range("J41").Value = 160
$ws.Range("J42").Value = 1025
$ws.Range("J44").Value = 180
$ws.Range("J51").Value = 295
$ws.Range("J52").Value = 595
$ws.Range("J53").Value = 342
$ws.Range("J54").Value = 453
$ws.Range("J55").Value = 344
$ws.Range("J60").Value = 139
$ws.Range("J63").Value = 85
$ws.Range("J65").Value = 585
$ws.Range("J66").Value = 70
$ws.Range("J67").Value = 894
$ws.Range("J68").Value = 51
$ws.Range("J70").Value = 34
$ws.Range("J75").Value = 71
$ws.Range("J79").Value = 668
$ws.Range("J83").Value = 473
$ws.Range("J84").Value = 197
$ws.Range("J85").Value = 990
$ws.Range("J88").Value = 249
$ws.Range("J91").Value = 272
$ws.Range("J92").Value = 76
$ws.Range("J94").Value = 250
$ws.Range("J95").Value = 339
$ws.Range("J97").Value = 211
$ws.Range("J99").Value = 364
$ws.Range("J101").Value = 23734

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 210
$ws.Range("J4").Value = 26

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 114
$ws.Range("J3").Value = 74
$ws.Range("J6").Value = 177
$ws.Range("J7").Value = 396

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 263
$ws.Range("J3").Value = 351
$ws.Range("J5").Value = 25
$ws.Range("J6").Value = 285
$ws.Range("J7").Value = 990

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 141
$ws.Range("J6").Value = 245
$ws.Range("J7").Value = 595

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 59
$ws.Range("J6").Value = 228
$ws.Range("J7").Value = 342

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 406
$ws.Range("J3").Value = 452
$ws.Range("J4").Value = 80
$ws.Range("J6").Value = 525
$ws.Range("J7").Value = 1500

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 127
$ws.Range("J7").Value = 473

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 254
$ws.Range("J6").Value = 379
$ws.Range("J7").Value = 1074

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 116
$ws.Range("J3").Value = 121
$ws.Range("J7").Value = 339

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 213
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 725

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 172
$ws.Range("J3").Value = 167
$ws.Range("J6").Value = 207
$ws.Range("J7").Value = 585

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 102
$ws.Range("J7").Value = 364

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J6").Value = 244
$ws.Range("J7").Value = 894

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 60
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 93
$ws.Range("J6").Value = 213
$ws.Range("J7").Value = 453

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 456
$ws.Range("J6").Value = 329
$ws.Range("J7").Value = 1302

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 205
$ws.Range("J6").Value = 271
$ws.Range("J7").Value = 701

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 58
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 218
$ws.Range("J6").Value = 547
$ws.Range("J7").Value = 1025

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 183
$ws.Range("J7").Value = 344

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 58
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 74
$ws.Range("J3").Value = 113
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 185
$ws.Range("J3").Value = 226
$ws.Range("J4").Value = 40
$ws.Range("J7").Value = 668

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 142
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 492

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 104
$ws.Range("J3").Value = 102
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 139
$ws.Range("J7").Value = 250

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 79
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J2").Value = 33
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 36
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 211

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 24
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J2").Value = 53
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 119
$ws.Range("J7").Value = 249

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value = 66
$ws.Range("J7").Value = 295

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 51
$ws.Range("J7").Value = 139

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J2").Value = 21
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 96
